$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23,81 TL - 23,81 TL"
$ws.Range("I2").Value = "18 TL - 18 TL"
$ws.Range("J2").Value = "15 TL - 15 TL"
$ws.Range("C3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("D3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("E3").Value = ""
$ws.Range("G3").Value = ""
$ws.Range("I3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("C4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("D4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("E4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("I4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("C5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("D5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("E5").Value = ""
$ws.Range("G5").Value = ""
$ws.Range("I5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("C6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("D6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("E6").Value = ""
$ws.Range("G6").Value = ""
$ws.Range("I6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("K6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("D7").Value = "%1,6"
$ws.Range("J7").Value = "%2,5"
$ws.Range("C8").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("D8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("E8").Value = ""
$ws.Range("G8").Value = ""
$ws.Range("I8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("C9").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("D9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("E9").Value = ""
$ws.Range("G9").Value = ""
$ws.Range("I9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("C10").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("D10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("E10").Value = ""
$ws.Range("G10").Value = ""
$ws.Range("I10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("C11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("D11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("E11").Value = ""
$ws.Range("G11").Value = ""
$ws.Range("I11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("K11").Value = "3,05 TL - 6,09 TL - 76,17 TL"
$ws.Range("C12").Value = "WU: 1.000,01 USD–9,51 USD"
$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"
$ws.Range("G12").Value = ""
$ws.Range("K12").Value = "WU: ,USD–; Diğer: 404,16 TL–3.403,42 TL"
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 9.999.999.999.999 TL"
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("E13").Value = ""
$ws.Range("I13").Value = "Hesaba: Asgari 1 TL | Azami 6,09 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"
$ws.Range("C14").Value = "40.000 TL - 1.904,76 TL"
$ws.Range("D14").Value = "2.300 TL - 9.500 TL"
$ws.Range("E14").Value = ""
$ws.Range("G14").Value = ""
$ws.Range("K14").Value = "914,14 TL - 4.265,98 TL"
$ws.Range("K15").Value = "%0,3 Asgari Tutar: 197,72 TL Azami Tutar: 197,72 TL / 249,13 TL"
$ws.Range("K17").Value = "%0,9 Asgari Tutar: 197,72 TL Azami Tutar: 197,72 TL / 2.528,89 TL"
$ws.Range("K20").Value = "122,59 TL"
$ws.Range("K21").Value = "%0,9 Asgari Tutar: 371,72 TL Azami Tutar: 371,72 TL / 2.022,72 TL"
$ws.Range("K22").Value = "%0,3 Asgari Tutar: 61,3 TL Azami Tutar: 61,3 TL / 7.596,55 TL"
$ws.Range("K23").Value = "54 TL"
$ws.Range("K24").Value = "371,72 TL"
$ws.Range("K25").Value = "312 TL"
